{"js": "// The document body contains a single paragraph used by M2Doc to report a\n// template/runtime version mismatch (an orange \"<---M2Doc version mismatch:\n// template is 3.1.1 and runtime is 3.2.0\" marker surrounded by padding\n// whitespace runs). The fix removes that diagnostic text, leaving the\n// paragraph with its original single empty run (<w:r><w:t/></w:r>).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst paragraph = paragraphs.items[0];\n// Replacing with an empty string clears all the runs' text/formatting and\n// collapses the paragraph back down to a single empty run, matching the\n// expected OOXML exactly.\nparagraph.insertText(\"\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# The document body's single paragraph carries an M2Doc diagnostic marker:\n# padding whitespace + an orange \"<---M2Doc version mismatch: template is\n# 3.1.1 and runtime is 3.2.0\" message + trailing padding whitespace, all\n# following an initial empty run. The fix strips that diagnostic text back\n# out, restoring the paragraph to just its original empty run.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n# Match the padding + marker text (including the leading four spaces and the\n# trailing four spaces) and replace it with nothing, leaving the paragraph's\n# leading empty run (<w:r><w:t/></w:r>) untouched.\n$find.Execute(\n    \"    <---M2Doc version mismatch: template is 3.1.1 and runtime is 3.2.0    \",\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"\",\n    2\n)\n"}
